# Edit script: insert two new weekly price rows for "Zapallo italiano"
# (Vega Modelo de Temuco) into the existing data table.
#
# Effect: rows currently at 213-305 shift down by two rows (to 215-307),
# and two brand-new rows of data are written into rows 213-214.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 213; this pushes the
# existing rows 213:305 down to 215:307, exactly matching the target layout.
$ws.Range("A213:A214").EntireRow.Insert()

# --- New row 213 ---
$ws.Cells.Item(213, 1).Value2  = 10
$ws.Cells.Item(213, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(213, 3).Value2  = "La Araucanía"
$ws.Cells.Item(213, 4).Value2  = 44523
$ws.Cells.Item(213, 5).Value2  = 9
$ws.Cells.Item(213, 6).Value2  = 100112032
$ws.Cells.Item(213, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(213, 8).Value2  = "Bola 8"
$ws.Cells.Item(213, 9).Value2  = "Primera"
$ws.Cells.Item(213, 10).Value2 = 40
$ws.Cells.Item(213, 11).Value2 = 8000
$ws.Cells.Item(213, 12).Value2 = 8000
$ws.Cells.Item(213, 13).Value2 = 8000
$ws.Cells.Item(213, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(213, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(213, 16).Value2 = 133
$ws.Cells.Item(213, 17).Value2 = 60
$ws.Cells.Item(213, 18).Value2 = "Hortaliza"

# --- New row 214 ---
$ws.Cells.Item(214, 1).Value2  = 10
$ws.Cells.Item(214, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(214, 3).Value2  = "La Araucanía"
$ws.Cells.Item(214, 4).Value2  = 44523
$ws.Cells.Item(214, 5).Value2  = 9
$ws.Cells.Item(214, 6).Value2  = 100112032
$ws.Cells.Item(214, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(214, 8).Value2  = "Sin especificar"
$ws.Cells.Item(214, 9).Value2  = "Primera"
$ws.Cells.Item(214, 10).Value2 = 300
$ws.Cells.Item(214, 11).Value2 = 8000
$ws.Cells.Item(214, 12).Value2 = 10000
$ws.Cells.Item(214, 13).Value2 = 9000
$ws.Cells.Item(214, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(214, 15).Value2 = "Región del Maule"
$ws.Cells.Item(214, 16).Value2 = 150
$ws.Cells.Item(214, 17).Value2 = 60
$ws.Cells.Item(214, 18).Value2 = "Hortaliza"
